# Fruta / hortaliza, semanal
# Insert this week's newest price readings (date 44474, Murcott / Mandarina)
# at the top of the data block, pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 126-127 (existing rows 126-143 shift to 128-145).
$ws.Range("A126:A127").EntireRow.Insert()

# New row 126: Murcott, Primera, Provincia de Limarí
$ws.Cells.Item(126, 1).Value = 7
$ws.Cells.Item(126, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(126, 3).Value = "Ñuble"
$ws.Cells.Item(126, 4).Value = 44474
$ws.Cells.Item(126, 5).Value = 16
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100102
$ws.Cells.Item(126, 8).Value = "Cítricos"
$ws.Cells.Item(126, 9).Value = 100102004
$ws.Cells.Item(126, 10).Value = "Mandarina"
$ws.Cells.Item(126, 11).Value = "Murcott"
$ws.Cells.Item(126, 12).Value = "Primera"
$ws.Cells.Item(126, 13).Value = 300
$ws.Cells.Item(126, 14).Value = 6000
$ws.Cells.Item(126, 15).Value = 6500
$ws.Cells.Item(126, 16).Value = 6250
$ws.Cells.Item(126, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(126, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(126, 19).Value = 625
$ws.Cells.Item(126, 20).Value = 10

# New row 127: Murcott, Segunda, Provincia de Limarí
$ws.Cells.Item(127, 1).Value = 7
$ws.Cells.Item(127, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(127, 3).Value = "Ñuble"
$ws.Cells.Item(127, 4).Value = 44474
$ws.Cells.Item(127, 5).Value = 16
$ws.Cells.Item(127, 6).Value = "Fruta"
$ws.Cells.Item(127, 7).Value = 100102
$ws.Cells.Item(127, 8).Value = "Cítricos"
$ws.Cells.Item(127, 9).Value = 100102004
$ws.Cells.Item(127, 10).Value = "Mandarina"
$ws.Cells.Item(127, 11).Value = "Murcott"
$ws.Cells.Item(127, 12).Value = "Segunda"
$ws.Cells.Item(127, 13).Value = 240
$ws.Cells.Item(127, 14).Value = 5000
$ws.Cells.Item(127, 15).Value = 5500
$ws.Cells.Item(127, 16).Value = 5250
$ws.Cells.Item(127, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(127, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(127, 19).Value = 525
$ws.Cells.Item(127, 20).Value = 10

# Make sure the date cells keep the workbook's date-time number format
# (the row insert should already have propagated this, but set explicitly
# to be safe).
$ws.Cells.Item(126, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(127, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "done"
